$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: update semantic dimension labels for "sexo" and "mes-y-ano" columns
# (iaest-measure: -> iaest-dimension:, since both are now treated as dimensions)
$ws.Range("G3").Value = "iaest-dimension:sexo"
$ws.Range("H3").Value = "iaest-dimension:mes-y-ano"

# Row 4: both columns become dimensions ("dim") instead of measures ("medida")
$ws.Range("G4").Value = "dim"
$ws.Range("H4").Value = "dim"

# Row 5: "sexo" now references a controlled vocabulary concept (skos:Concept)
# instead of a plain string
$ws.Range("G5").Value = "skos:Concept"

# Row 6 (new row): mapping file reference for the "sexo" dimension values
$ws.Range("G6").Value = "mapping-sexo.xlsx"

# Match the formatting used by the rest of the column (reuse the same cell
# style as G5 rather than leaving G6 with the workbook default style)
$ws.Range("G5").Copy()
$ws.Range("G6").PasteSpecial(-4122)
$excel.CutCopyMode = $false
